$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before BG (shifts old BG -> BH, old BH -> BI).
$ws.Columns("BG").Insert()

# New header cell for the freshly-scraped timestamp column.
$ws.Range("BG1").Value2 = "2026-01-30 09:27:43"

# Carry forward the latest known price into the new column for every
# product row that still had a numeric price in BF (rows 2-80).
# Rows 81-206 had already stopped updating (empty BF), so BG stays empty.
for ($r = 2; $r -le 80; $r++) {
    $ws.Cells.Item($r, 59).Value2 = $ws.Cells.Item($r, 58).Value2
}
